# Weekly update: two new daily price records were added to the
# "Femacal de La Calera - Albahaca" sheet. Excel keeps the sheet sorted by
# date descending-ish (actually just by arrival order), so the two new
# rows are inserted in the middle of the existing data, pushing the
# rows that follow down.
#
# New record #1 is inserted at row 20 (pushes old row 20.. down by one).
# New record #2 is inserted at row 86 of the resulting sheet (pushes the
# remaining old rows down by one more).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new record at row 20 -----------------------------------------
$ws.Rows.Item(20).Insert()

$ws.Range("A20").Value = 3
$ws.Range("B20").Value = 'Femacal de La Calera'
$ws.Range("C20").Value = 'Coquimbo'
$ws.Range("D20").Value = 44635
$ws.Range("E20").Value = 5
$ws.Range("F20").Value = 100112052
$ws.Range("G20").Value = 'Albahaca'
$ws.Range("H20").Value = 'Sin especificar'
$ws.Range("I20").Value = 'Primera'
$ws.Range("J20").Value = 120
$ws.Range("K20").Value = 5000
$ws.Range("L20").Value = 5500
$ws.Range("M20").Value = 5250
$ws.Range("N20").Value = '$/docena de matas'
$ws.Range("O20").Value = 'Provincia de Quillota'
$ws.Range("P20").Value = 875
$ws.Range("Q20").Value = 6
$ws.Range("R20").Value = 'Hortaliza'

# --- Insert new record at row 86 (post first insert numbering) ----------
$ws.Rows.Item(86).Insert()

$ws.Range("A86").Value = 3
$ws.Range("B86").Value = 'Femacal de La Calera'
$ws.Range("C86").Value = 'Coquimbo'
$ws.Range("D86").Value = 44634
$ws.Range("E86").Value = 5
$ws.Range("F86").Value = 100112052
$ws.Range("G86").Value = 'Albahaca'
$ws.Range("H86").Value = 'Sin especificar'
$ws.Range("I86").Value = 'Primera'
$ws.Range("J86").Value = 105
$ws.Range("K86").Value = 5000
$ws.Range("L86").Value = 5500
$ws.Range("M86").Value = 5262
$ws.Range("N86").Value = '$/docena de matas'
$ws.Range("O86").Value = 'Provincia de Quillota'
$ws.Range("P86").Value = 877
$ws.Range("Q86").Value = 6
$ws.Range("R86").Value = 'Hortaliza'
